$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The googleMapsClient script was re-appending the same geocoded
# rows on every run instead of overwriting them, so rows 2-14 are
# duplicate/garbage data -- remove them, leaving only the header row.
$ws.Range("A2:C14").EntireRow.Delete() | Out-Null

# The buggy run also wrote blank/garbage into the "Address" header
# cell -- clear it out so row 1 only has the Latitude/Longitude
# headers left.
$ws.Range("A1").ClearContents() | Out-Null

# Touch A1's formatting (a no-op against the sheet's default font) so
# column A stays part of the sheet's used range/dimension even though
# its value was cleared, matching the header row's original layout.
$ws.Range("A1").Font.Bold = $false
